# Applies the "Extent report code" commit:
#   - Supervisor re-assigned from Aravind ET / Aravinda ET -> Palak Garg
#   - OrgUnit path changed from the old Singapore/Mangalore paths to
#     India>South>Bangalore>ProductQA (short name "ProductQA") on the
#     Create / Edit / Delete sheets
#   - assorted sheet-view / selection bookkeeping that goes along with it

$wb = $excel.ActiveWorkbook

$newOrgPath  = "India>South>Bangalore>ProductQA"
$newOrgShort = "ProductQA"
$newSupervisor = "Palak Garg"

# ---- Create sheet ----------------------------------------------------
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("E2:E3").Value2 = $newOrgPath
$wsCreate.Range("G2:G3").Value2 = $newSupervisor
$wsCreate.Range("I2:I3").Value2 = $newOrgShort
$wsCreate.Columns.Item(9).ColumnWidth = 12.7369791666667

# ---- Edit sheet --------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("E2:E3").Value2 = $newOrgPath
$wsEdit.Range("G2:G3").Value2 = $newSupervisor
$wsEdit.Range("I2:I3").Value2 = $newOrgShort
# D2 carried a quote-prefixed text style; rewriting the value drops it
# so D2 matches the plain (non quote-prefixed) style already used by D3.
$wsEdit.Range("D2").Value2 = $wsEdit.Range("D2").Value2

# ---- Delete sheet --------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("E2:E3").Value2 = $newOrgPath
$wsDelete.Range("G2:G3").Value2 = $newSupervisor
$wsDelete.Range("I2:I3").Value2 = $newOrgShort

# ---- view bookkeeping (selection / active sheet) ----------------------
# Select on the non-active sheets first; whichever sheet is selected last
# ends up the workbook's active tab, so "Create" is activated last below.
$wsEdit.Range("E2").Select() | Out-Null
$wsDelete.Range("I5").Select() | Out-Null

$wsCreate.Activate() | Out-Null
$wsCreate.Range("B4").Select() | Out-Null
